# JobHistory/Read - write the suite row into the Test Suite Statistics sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New suite row: Read - 0 automated out of 7 total test cases, status "Testing"
# (mirrors the existing ListView row in row 2).
$ws.Range("A3").Value = "Read"
$ws.Range("B3").Value = 0
$ws.Range("C3").Value = 7
$ws.Range("D3").Value = "Testing"

# Move the active selection to D4, matching the author's saved cursor position.
$ws.Range("D4").Select() | Out-Null
